$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order matches: C,F,G,H,I,K,L,M,N,O,P,Q,R,S,T,U,V
$columns = @("C","F","G","H","I","K","L","M","N","O","P","Q","R","S","T","U","V")

$data = @(
    @(4, 325, "bedrooms", "bedrooms", "target", "j", "stimuli/img_2js6m.png", 40.02777777777778, 20.88888888888889, 30.45833333333334, 36, 2, 2, 2, 2, 2, 2),
    @(4, 326, "bedrooms", "bedrooms", "target", "j", "stimuli/img_cmyvx.png", 64.25, 40.09375, 52.171875, 32, 4, 4, 4, 4, 4, 4),
    @(4, 327, "bedrooms", "bedrooms", "target", "j", "stimuli/img_juob3.png", 79.92105263157895, 59.78947368421053, 69.85526315789474, 38, 7, 7, 7, 7, 7, 7),
    @(4, 328, "bedrooms", "living_rooms", "distractor", "f", "stimuli/img_vh7v8.png", 78.70454545454545, 59.63636363636363, 69.17045454545455, 44, 7, 7, 7, 7, 7, 7),
    @(4, 329, "bedrooms", "bedrooms", "target", "j", "stimuli/img_3h4c9.png", 85.47619047619048, 67.26190476190476, 76.36904761904762, 42, 9, 9, 9, 9, 9, 9),
    @(4, 330, "bedrooms", "bedrooms", "target", "j", "stimuli/img_wyctg.png", 33.44736842105263, 11.39473684210526, 22.42105263157895, 38, 1, 1, 1, 1, 1, 1),
    @(4, 331, "bedrooms", "kitchens", "distractor", "f", "stimuli/img_kugyw.png", 74.25, 54.10714285714285, 64.17857142857143, 28, 6, 6, 6, 6, 6, 6),
    @(4, 332, "bedrooms", "bedrooms", "target", "j", "stimuli/img_oou46.png", 75.70270270270271, 54.86486486486486, 65.28378378378379, 37, 6, 6, 6, 6, 6, 6),
    @(4, 333, "bedrooms", "kitchens", "distractor", "f", "stimuli/img_xguy9.png", 78.21621621621621, 57.24324324324324, 67.72972972972973, 37, 7, 7, 7, 7, 7, 7),
    @(4, 334, "bedrooms", "bedrooms", "target", "j", "stimuli/img_rvssl.png", 74.25, 54.33333333333334, 64.29166666666667, 36, 6, 6, 6, 6, 6, 6),
    @(4, 335, "bedrooms", "living_rooms", "distractor", "f", "stimuli/img_89dvt.png", 81.09756097560975, 64.6829268292683, 72.89024390243902, 41, 8, 8, 8, 8, 8, 8),
    @(4, 336, "bedrooms", "living_rooms", "distractor", "f", "stimuli/img_eh0no.png", 53.66666666666666, 36.02564102564103, 44.84615384615385, 39, 3, 3, 3, 4, 3, 4),
    @(4, 337, "bedrooms", "living_rooms", "distractor", "f", "stimuli/img_7lz7m.png", 51.5531914893617, 32.87234042553192, 42.21276595744681, 47, 3, 3, 3, 3, 3, 3),
    @(4, 338, "bedrooms", "bedrooms", "target", "j", "stimuli/img_72fmj.png", 53.87179487179487, 36.02564102564103, 44.94871794871795, 39, 3, 3, 3, 3, 3, 3),
    @(4, 339, "bedrooms", "living_rooms", "distractor", "f", "stimuli/img_24rt2.png", 55.26829268292683, 34.19512195121951, 44.73170731707317, 41, 3, 3, 3, 4, 4, 3),
    @(4, 340, "bedrooms", "bedrooms", "target", "j", "stimuli/img_f4jxo.png", 82.91666666666667, 65.52777777777777, 74.22222222222223, 36, 8, 8, 8, 8, 8, 8),
    @(4, 341, "bedrooms", "kitchens", "distractor", "f", "stimuli/img_lpj57.png", 74.77777777777777, 54.44444444444444, 64.61111111111111, 27, 6, 6, 6, 6, 6, 6),
    @(4, 342, "bedrooms", "kitchens", "distractor", "f", "stimuli/img_qihxi.png", 76.72222222222223, 56.33333333333334, 66.52777777777779, 36, 7, 7, 7, 7, 7, 7),
    @(4, 343, "bedrooms", "kitchens", "distractor", "f", "stimuli/img_eppte.png", 78.42424242424242, 57.03030303030303, 67.72727272727272, 33, 7, 7, 7, 7, 7, 7),
    @(4, 344, "bedrooms", "bedrooms", "target", "j", "stimuli/img_1vq1v.png", 69.42857142857143, 46.59523809523809, 58.01190476190476, 42, 5, 5, 5, 5, 5, 5),
    @(4, 345, "bedrooms", "bedrooms", "target", "j", "stimuli/img_uxxo0.png", 71.74418604651163, 48.44186046511628, 60.09302325581395, 43, 5, 5, 5, 5, 5, 5),
    @(4, 346, "bedrooms", "kitchens", "distractor", "f", "stimuli/img_cxpff.png", 74.92307692307692, 53.28205128205128, 64.1025641025641, 39, 6, 6, 6, 6, 6, 6),
    @(4, 347, "bedrooms", "kitchens", "distractor", "f", "stimuli/img_60242.png", 78.33333333333333, 57.57575757575758, 67.95454545454545, 33, 7, 7, 7, 7, 7, 7),
    @(4, 348, "bedrooms", "bedrooms", "target", "j", "stimuli/img_t2ioc.png", 88.18918918918919, 74.05405405405405, 81.12162162162161, 37, 10, 10, 10, 10, 10, 10),
    @(4, 349, "bedrooms", "bedrooms", "target", "j", "stimuli/img_x0u5z.png", 92, 78.16216216216216, 85.08108108108108, 37, 10, 10, 10, 10, 10, 10),
    @(4, 350, "bedrooms", "living_rooms", "distractor", "f", "stimuli/img_3jnt7.png", 49.52272727272727, 35.25, 42.38636363636364, 44, 3, 3, 3, 3, 3, 4),
    @(4, 351, "bedrooms", "bedrooms", "target", "j", "stimuli/img_a9acb.png", 77.11428571428571, 58.42857142857143, 67.77142857142857, 35, 7, 7, 7, 7, 7, 7),
    @(4, 352, "bedrooms", "kitchens", "distractor", "f", "stimuli/img_79b5l.png", 72.74285714285715, 53.31428571428572, 63.02857142857143, 35, 6, 6, 6, 6, 6, 6),
    @(4, 353, "bedrooms", "bedrooms", "target", "j", "stimuli/img_jp28n.png", 65.02564102564102, 44.97435897435897, 55, 39, 4, 4, 4, 4, 4, 5),
    @(4, 354, "bedrooms", "living_rooms", "distractor", "f", "stimuli/img_cehin.png", 78.86363636363636, 60.02272727272727, 69.44318181818181, 44, 7, 7, 7, 7, 7, 7),
    @(4, 355, "bedrooms", "living_rooms", "distractor", "f", "stimuli/img_xpco9.png", 81.55555555555556, 64.68888888888888, 73.12222222222222, 45, 8, 8, 8, 8, 8, 8),
    @(4, 356, "bedrooms", "kitchens", "distractor", "f", "stimuli/img_0mhms.png", 78, 55.68571428571428, 66.84285714285714, 35, 7, 7, 7, 7, 7, 7),
    @(4, 357, "bedrooms", "bedrooms", "target", "j", "stimuli/img_le8uf.png", 12.88888888888889, 9.222222222222221, 11.05555555555556, 36, 1, 1, 1, 1, 1, 1),
    @(4, 358, "bedrooms", "living_rooms", "distractor", "f", "stimuli/img_emh91.png", 82.06666666666666, 63.33333333333334, 72.7, 45, 8, 8, 8, 8, 8, 8),
    @(4, 359, "bedrooms", "bedrooms", "target", "j", "stimuli/img_e26ut.png", 81.07692307692308, 61.28205128205128, 71.17948717948718, 39, 8, 8, 8, 8, 8, 8),
    @(4, 360, "bedrooms", "kitchens", "distractor", "f", "stimuli/img_e0hwx.png", 78.12121212121212, 55.36363636363637, 66.74242424242425, 33, 7, 7, 7, 7, 7, 7),
    @(4, 361, "bedrooms", "bedrooms", "target", "j", "stimuli/img_5il0t.png", 48.09523809523809, 30.90476190476191, 39.5, 42, 2, 2, 2, 2, 2, 2),
    @(4, 362, "bedrooms", "living_rooms", "distractor", "f", "stimuli/img_3m61b.png", 81.97619047619048, 63.23809523809524, 72.60714285714286, 42, 8, 8, 8, 8, 8, 8),
    @(4, 363, "bedrooms", "bedrooms", "target", "j", "stimuli/img_cogrz.png", 60.5, 39.71428571428572, 50.10714285714286, 42, 3, 3, 3, 3, 3, 3),
    @(4, 364, "bedrooms", "bedrooms", "target", "j", "stimuli/img_fnu4h.png", 85.87179487179488, 70.71794871794872, 78.2948717948718, 39, 9, 9, 9, 9, 9, 9)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowNum = $r + 2
    $rowData = $data[$r]
    for ($c = 0; $c -lt $columns.Count; $c++) {
        $cellRef = "$($columns[$c])$rowNum"
        $ws.Range($cellRef).Value = $rowData[$c]
    }
}
